# feat: add front image (fixes #2)
#
# Adds a third column ("image") to the sheet that records which image
# file illustrates each row, and updates the print orientation /
# selection the way Excel does when a picture is inserted and the
# workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C: header + one filename per data row.
# Written bottom-up so the shared-string table is populated in the same
# order ("simple stick figure.png", "umbrella.png", "image") as the
# authoritative export.
$ws.Range("C3").Value = "simple stick figure.png"
$ws.Range("C2").Value = "umbrella.png"
$ws.Range("C1").Value = "image"

# Excel switches the page to portrait orientation once a picture/print
# setup is touched.
$ws.PageSetup.Orientation = 1

# Final selection left on the sheet after the edit.
$ws.Range("E9").Select() | Out-Null
